$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1: top+bottom thin border (matches existing borderId 4)
$c = $ws1.Range("C1")
$c.Borders.LineStyle = 1
$c.Borders.Weight = 2
$c.Borders.Item(7).LineStyle = 0
$c.Borders.Item(10).LineStyle = 0

# D1: top+bottom+right thin border (matches existing borderId 5)
$d = $ws1.Range("D1")
$d.Borders.LineStyle = 1
$d.Borders.Weight = 2
$d.Borders.Item(7).LineStyle = 0

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

$c2 = $ws2.Range("C1")
$c2.Borders.LineStyle = 1
$c2.Borders.Weight = 2
$c2.Borders.Item(7).LineStyle = 0
$c2.Borders.Item(10).LineStyle = 0

$d2 = $ws2.Range("D1")
$d2.Borders.LineStyle = 1
$d2.Borders.Weight = 2
$d2.Borders.Item(7).LineStyle = 0

$f2 = $ws2.Range("F1")
$f2.Borders.LineStyle = 1
$f2.Borders.Weight = 2
$f2.Borders.Item(7).LineStyle = 0
$f2.Borders.Item(10).LineStyle = 0

$g2 = $ws2.Range("G1")
$g2.Borders.LineStyle = 1
$g2.Borders.Weight = 2
$g2.Borders.Item(7).LineStyle = 0

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
